# GlobalHIVMerged.xlsx - strip the trailing " (Source: ...)" annotation from
# the "Descrimination Percent" column (C2:C60) so the cells hold plain
# numeric percentages instead of text, using a regex find/replace
# (matches the commit: find ' *\(.*\)' / replace with '').

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$dataRange = $ws.Range("C2:C60")

foreach ($cell in $dataRange.Cells) {
    $raw = $cell.Value()
    if ($raw -ne $null) {
        $stripped = [regex]::Replace([string]$raw, ' *\(.*\)\s*$', '')
        $cell.Value = [double]$stripped
    }
}

# Restore the active selection recorded in the saved file.
$ws.Range("I16").Select()
